$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the prices in column D for rows 31-33
$ws.Range("D31").Value = 3985.166
$ws.Range("D32").Value = 1646.582
$ws.Range("D33").Value = 3379.108

# Touch (re-merge) the merged ranges so their order in the saved
# workbook matches the order produced by the original edit.
$ws.Range("B32:C32").UnMerge()
$ws.Range("B32:C32").Merge()
$ws.Range("A10:D10").UnMerge()
$ws.Range("A10:D10").Merge()
$ws.Range("B33:C33").UnMerge()
$ws.Range("B33:C33").Merge()
$ws.Range("A11:D11").UnMerge()
$ws.Range("A11:D11").Merge()
$ws.Range("B31:C31").UnMerge()
$ws.Range("B31:C31").Merge()
$ws.Range("A1:D1").UnMerge()
$ws.Range("A1:D1").Merge()
$ws.Range("B30:C30").UnMerge()
$ws.Range("B30:C30").Merge()
$ws.Range("A9:D9").UnMerge()
$ws.Range("A9:D9").Merge()
